# Weekly price-sheet refresh for "Bruselas (repollito)" / Vega Modelo de Temuco.
# Rows 14-35 already contain data (dates/volumes/prices shuffled around as new
# weekly observations land); rows 36-41 are brand-new observations appended at
# the bottom of the table. Column D carries the "YYYY-MM-DD HH:MM:SS" date
# style (s="2" in the original sheet); every other populated column is plain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, date-serial(D), volumen(J), precioMin(K), precioMax(L), precioProm(M), precioKg(P)
$data = @(
    @(14, 44434, 55, 25000, 25000, 25000, 2500),
    @(15, 44428, 30, 25000, 25000, 25000, 2500),
    @(16, 44384, 40, 25000, 25000, 25000, 2500),
    @(17, 44349, 45, 24000, 24000, 24000, 2400),
    @(18, 44385, 80, 25000, 25000, 25000, 2500),
    @(19, 44427, 40, 25000, 25000, 25000, 2500),
    @(20, 44413, 40, 25000, 25000, 25000, 2500),
    @(21, 44421, 55, 25000, 25000, 25000, 2500),
    @(22, 44400, 12, 24000, 24000, 24000, 2400),
    @(23, 44426, 30, 25000, 25000, 25000, 2500),
    @(24, 44390, 15, 25000, 25000, 25000, 2500),
    @(25, 44354, 30, 24000, 24000, 24000, 2400),
    @(26, 44410, 50, 25000, 25000, 25000, 2500),
    @(27, 44412, 50, 25000, 25000, 25000, 2500),
    @(28, 44371, 50, 25000, 25000, 25000, 2500),
    @(29, 44435, 185, 25000, 27000, 25162, 2516),
    @(30, 44431, 65, 25000, 25000, 25000, 2500),
    @(31, 44405, 40, 25000, 25000, 25000, 2500),
    @(32, 44350, 40, 24000, 25000, 24375, 2438),
    @(33, 44389, 65, 25000, 25000, 25000, 2500),
    @(34, 44417, 15, 25000, 25000, 25000, 2500),
    @(35, 44419, 25, 25000, 25000, 25000, 2500),
    @(36, 44420, 55, 25000, 25000, 25000, 2500),
    @(37, 44433, 25, 25000, 25000, 25000, 2500),
    @(38, 44382, 50, 25000, 25000, 25000, 2500),
    @(39, 44355, 25, 23000, 24000, 23400, 2340),
    @(40, 44376, 45, 23000, 23000, 23000, 2300),
    @(41, 44432, 15, 27000, 27000, 27000, 2700)
)

# Columns that are identical across every data row of this subset.
$constCols = @{
    1  = 10;                          # A Mercado ID
    2  = "Vega Modelo de Temuco";     # B Mercado
    3  = "La Araucanía";              # C Región
    5  = 9;                           # E Codreg
    6  = 100112035;                   # F Categoría ID
    7  = "Bruselas (repollito)";      # G Categoría
    8  = "Sin especificar";           # H Variedad
    9  = "Primera";                   # I Calidad
    14 = "`$/malla 10 kilos";         # N Unidad de comercialización
    15 = "Provincia de Quillota";     # O Origen
    17 = 10;                          # Q Kg o Unidades
    18 = "Hortaliza"                  # R Clasificación
}

$dateFormat = $ws.Cells.Item(35, 4).NumberFormat

foreach ($row in $data) {
    $r = $row[0]

    # Make sure every "always the same" column is populated (needed for the
    # brand-new rows 36-41; harmless no-op for the pre-existing rows).
    foreach ($col in $constCols.Keys) {
        $ws.Cells.Item($r, $col).Value = $constCols[$col]
    }

    $ws.Cells.Item($r, 4).Value = $row[1]            # D Fecha
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 10).Value = $row[2]           # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]           # K Precio mínimo
    $ws.Cells.Item($r, 12).Value = $row[4]           # L Precio máximo
    $ws.Cells.Item($r, 13).Value = $row[5]           # M Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]           # P Precio $/Kg
}
